$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.8
$ws.Range("I2").Value = 2.1
$ws.Range("J2").Value = 4.5
$ws.Range("L2").Value = 2.88
$ws.Range("AA2").Value = 41
$ws.Range("AH2").Value = 8.5
$ws.Range("AJ2").Value = 19
$ws.Range("AO2").Value = 23
$ws.Range("AX2").Value = 12

# Row 3
$ws.Range("G3").Value = 1.85
$ws.Range("I3").Value = 5
$ws.Range("L3").Value = 6
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 4.75
$ws.Range("AX3").Value = 34

# Row 4
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67

# Row 5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10

# Row 6
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1.73
$ws.Range("R6").Value = 2.08

# Row 11
$ws.Range("G11").Value = 2.1
$ws.Range("I11").Value = 3.4
$ws.Range("N11").Value = 7.5
$ws.Range("W11").Value = 6.5
$ws.Range("AA11").Value = 21
$ws.Range("AD11").Value = 6.5
$ws.Range("AE11").Value = 19
$ws.Range("AG11").Value = 8.5
$ws.Range("AO11").Value = 13
$ws.Range("AZ11").Value = 67
